# Starlendar: add a "Dozenal" calendar sheet alongside the existing
# (now renamed) "Decimal" sheet, drop the unused third sheet, and make
# the new Dozenal sheet the active one.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Rename Sheet1 -> Decimal (content/layout is unchanged there).
$wb.Worksheets.Item("Sheet1").Name = "Decimal"
$decimal = $wb.Worksheets.Item("Decimal")

# 2. Build the Dozenal sheet as a copy of Decimal (same layout, styles,
#    merges, column widths, etc.) placed right after the old Sheet2.
$sheet2 = $wb.Worksheets.Item("Sheet2")
$decimal.Copy($null, $sheet2)

# 3. Drop the now-redundant blank placeholder sheets.
[void]$sheet2.Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()

# 4. Name the freshly copied sheet and make it the active tab.
$dozenal = $wb.Worksheets.Item("Decimal (2)")
$dozenal.Name = "Dozenal"
$dozenal.Activate()
[void]$dozenal.Range("V22").Select()

# 5. Rewrite the day numbers for the dozenal (base-12) week layout: each
#    5-day block before was 1,7,13,19,25 (decimal) and becomes
#    1,7,11,17,21 (dozenal), with "X" (10) and "E" (11) replacing the
#    two-digit decimal numbers 10/11 and 22/23 respectively.
#    Cells are touched in the same order the new shared strings first
#    appear ("X", "E", "1X", "1E", then the "Ultra days" header label)
#    so newly-minted shared-string indices land the same way.
$blocks = @("B","C","D","E","F"), ("G","H","I","J","K"), ("L","M","N","O","P")

foreach ($cols in $blocks) {
    $dozenal.Range($cols[1] + "10").Value = "X"
}
foreach ($cols in $blocks) {
    $dozenal.Range($cols[1] + "11").Value = "E"
}
foreach ($cols in $blocks) {
    $dozenal.Range($cols[3] + "10").Value = "1X"
}
foreach ($cols in $blocks) {
    $dozenal.Range($cols[3] + "11").Value = "1E"
}

# 6. Re-point the header label in R1 (was "Leapday (U6)" on Decimal) to
#    the dozenal-specific "Ultra days" label.
$dozenal.Range("R1").Value = "Ultra days"

foreach ($cols in $blocks) {
    $dozenal.Range($cols[1] + "7").Value  = 7
    $dozenal.Range($cols[2] + "7").Value  = 11
    $dozenal.Range($cols[3] + "7").Value  = 17
    $dozenal.Range($cols[4] + "7").Value  = 21

    $dozenal.Range($cols[1] + "8").Value  = 8
    $dozenal.Range($cols[2] + "8").Value  = 12
    $dozenal.Range($cols[3] + "8").Value  = 18
    $dozenal.Range($cols[4] + "8").Value  = 22

    $dozenal.Range($cols[1] + "9").Value  = 9
    $dozenal.Range($cols[2] + "9").Value  = 13
    $dozenal.Range($cols[3] + "9").Value  = 19
    $dozenal.Range($cols[4] + "9").Value  = 23

    $dozenal.Range($cols[2] + "10").Value = 14
    $dozenal.Range($cols[4] + "10").Value = 24

    $dozenal.Range($cols[2] + "11").Value = 15
    $dozenal.Range($cols[4] + "11").Value = 25

    $dozenal.Range($cols[0] + "12").Value = 6
    $dozenal.Range($cols[1] + "12").Value = 10
    $dozenal.Range($cols[2] + "12").Value = 16
    $dozenal.Range($cols[3] + "12").Value = 20
    $dozenal.Range($cols[4] + "12").Value = 26
}
